$wb = $excel.ActiveWorkbook

# The daily HTTPS export data lives on the "Chart" sheet (sheet1.xml).
$chart = $wb.Worksheets.Item("Chart")

# Append the next day's row (2025-12-29) right after the last existing
# data row (2025-12-28), carrying forward the same "Pages" total (28)
# and "Non-HTTPS URLs" value (0) as the previous day.
$dateCell = $chart.Cells.Item(85, 1)

# Force the cell to be treated as literal text rather than letting the
# date-like string auto-convert into a date serial number, then restore
# the default (General) appearance so the cell matches its neighbours.
$dateCell.NumberFormat = "@"
$dateCell.Value = "2025-12-29"
$dateCell.ClearFormats()

$chart.Cells.Item(85, 2).Value = 0
$chart.Cells.Item(85, 3).Value = 28
